# Add two new rows (104, 105) of data to each of the 4 worksheets,
# continuing the existing date/remn_amt series.
#
# Sheet 1 (신화인터텍):  104 -> 45967 / 643   105 -> 45968 / 0
# Sheet 2 (드림시큐리티): 104 -> 45967 / 7403  105 -> 45968 / 0
# Sheet 3 (대원미디어):  104 -> 45967 / 2806  105 -> 45968 / 0
# Sheet 4 (성호전자):   104 -> 45967 / 1354  105 -> 45968 / 0

$wb = $excel.ActiveWorkbook

$newRows = @(
    @{ Row = 104; Date = 45967; Values = @(643, 7403, 2806, 1354) },
    @{ Row = 105; Date = 45968; Values = @(0, 0, 0, 0) }
)

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    foreach ($entry in $newRows) {
        $r = $entry.Row
        $dateCell = $ws.Cells.Item($r, 1)
        $amtCell = $ws.Cells.Item($r, 2)

        # Match the number format / style used by the date column above
        # (the last existing row, r-1) so the new cells render the same way.
        $dateCell.NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat

        $dateCell.Value = $entry.Date
        $amtCell.Value = $entry.Values[$i - 1]
    }
}
